# Update "想去人数" (interested-count) values in the "展览" and "全部类型"
# sheets, as per the commit that regenerated the gh-pages data output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    7  = 2755
    9  = 1767
    12 = 616
    15 = 151
    17 = 87
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
